$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two added columns
$ws.Range("G1").Value = "TransactionSpeedNS"
$ws.Range("H1").Value = "TransactionSpeedMS"

# Populate TransactionSpeedNS (ResolveTime - AttemptTime) and
# TransactionSpeedMS (rounded to nearest millisecond) for each data row
for ($r = 2; $r -le 51; $r++) {
    $attempt = $ws.Cells.Item($r, 3).Value2
    $resolve = $ws.Cells.Item($r, 4).Value2
    $ns = $resolve - $attempt
    $ms = [System.Math]::Round($ns / 1000000)

    $ws.Cells.Item($r, 7).Value = $ns
    $ws.Cells.Item($r, 8).Value = $ms
}
